$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 726.5
$ws.Range("I15").Value = 726.5
$ws.Range("K15").Value = 2179.5
$ws.Range("M15").Value = -2010.5
$ws.Range("H132").Value = 1825.758
$ws.Range("I132").Value = 1763.2833
$ws.Range("J132").Value = 3700
$ws.Range("K132").Value = 5289.8499
$ws.Range("L132").Value = 11100
$ws.Range("M132").Value = -2759.8499
$ws.Range("N132").Value = -16160
$ws.Range("H137").Value = 20367.166
$ws.Range("I137").Value = 12500
$ws.Range("J137").Value = 24300.75
$ws.Range("K137").Value = 37500
$ws.Range("L137").Value = 72902.25
$ws.Range("M137").Value = -34950
$ws.Range("N137").Value = -78002.25
$ws.Range("H141").Value = 1017.1795
$ws.Range("I141").Value = 873.8889
$ws.Range("J141").Value = 2736.6667
$ws.Range("K141").Value = 2621.6667
$ws.Range("L141").Value = 8210.000100000001
$ws.Range("M141").Value = 2558.3333
$ws.Range("N141").Value = -18570.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4681.508
$ws.Range("I32").Value = 4811.8037
$ws.Range("J32").Value = 4127.75
$ws.Range("K32").Value = 4811.8037
$ws.Range("L32").Value = 4127.75
$ws.Range("M32").Value = -4524.8037
$ws.Range("N32").Value = -4701.75
$ws.Range("H61").Value = 3117.6943
$ws.Range("I61").Value = 2426.08
$ws.Range("J61").Value = 4689.5454
$ws.Range("K61").Value = 2426.08
$ws.Range("L61").Value = 4689.5454
$ws.Range("M61").Value = -2214.08
$ws.Range("N61").Value = -5113.5454
$ws.Range("H74").Value = 24845.51
$ws.Range("J74").Value = 2348.125
$ws.Range("L74").Value = 2348.125
$ws.Range("N74").Value = -4096.125
$ws.Range("H77").Value = 24845.51
$ws.Range("J77").Value = 2348.125
$ws.Range("L77").Value = 11740.625
$ws.Range("N77").Value = -20476.625
$ws.Range("H102").Value = 7910
$ws.Range("J102").Value = 21331.334
$ws.Range("L102").Value = 21331.334
$ws.Range("N102").Value = -24575.334
$ws.Range("H132").Value = 85659.75999999999
$ws.Range("I132").Value = 2649.5
$ws.Range("K132").Value = 7948.5
$ws.Range("M132").Value = -5418.5
$ws.Range("H136").Value = 3117.6943
$ws.Range("I136").Value = 2426.08
$ws.Range("J136").Value = 4689.5454
$ws.Range("K136").Value = 7278.24
$ws.Range("L136").Value = 14068.6362
$ws.Range("M136").Value = -4728.24
$ws.Range("N136").Value = -19168.6362

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 52626.5
$ws.Range("J60").Value = 52626.5
$ws.Range("L60").Value = 52626.5
$ws.Range("N60").Value = -53824.5
$ws.Range("H96").Value = 77674.125
$ws.Range("I96").Value = 6348.5
$ws.Range("K96").Value = 6348.5
$ws.Range("M96").Value = -3602.5
$ws.Range("H134").Value = 3523.7727
$ws.Range("I134").Value = 2148.7646
$ws.Range("J134").Value = 8198.799999999999
$ws.Range("K134").Value = 6446.293799999999
$ws.Range("L134").Value = 24596.4
$ws.Range("M134").Value = -3911.293799999999
$ws.Range("N134").Value = -29666.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4580.4375
$ws.Range("I31").Value = 2671
$ws.Range("J31").Value = 7253.65
$ws.Range("K31").Value = 2671
$ws.Range("L31").Value = 7253.65
$ws.Range("M31").Value = -2376
$ws.Range("N31").Value = -7843.65
$ws.Range("H34").Value = 4580.4375
$ws.Range("I34").Value = 2671
$ws.Range("J34").Value = 7253.65
$ws.Range("K34").Value = 2671
$ws.Range("L34").Value = 7253.65
$ws.Range("M34").Value = -2469
$ws.Range("N34").Value = -7657.65
$ws.Range("H58").Value = 3091.4
$ws.Range("I58").Value = 2448.0588
$ws.Range("J58").Value = 4458.5
$ws.Range("K58").Value = 2448.0588
$ws.Range("L58").Value = 4458.5
$ws.Range("M58").Value = -2245.0588
$ws.Range("N58").Value = -4864.5
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H132").Value = 3064.3333
$ws.Range("I132").Value = 3068.9285
$ws.Range("K132").Value = 9206.7855
$ws.Range("M132").Value = -6676.7855
$ws.Range("H134").Value = 2489.5625
$ws.Range("I134").Value = 1888.8148
$ws.Range("K134").Value = 5666.4444
$ws.Range("M134").Value = -3131.4444
$ws.Range("H136").Value = 3091.4
$ws.Range("I136").Value = 2448.0588
$ws.Range("J136").Value = 4458.5
$ws.Range("K136").Value = 7344.176399999999
$ws.Range("L136").Value = 13375.5
$ws.Range("M136").Value = -4794.176399999999
$ws.Range("N136").Value = -18475.5
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3133
$ws.Range("J39").Value = 4999
$ws.Range("L39").Value = 14997
$ws.Range("N39").Value = -15585
$ws.Range("H122").Value = 1679.7894
$ws.Range("J122").Value = 1701
$ws.Range("L122").Value = 15309
$ws.Range("N122").Value = -20209
$ws.Range("H129").Value = 1429.16
$ws.Range("I129").Value = 723.2857
$ws.Range("K129").Value = 2169.8571
$ws.Range("M129").Value = 2830.1429
$ws.Range("H131").Value = 1576.6666
$ws.Range("J131").Value = 1710.7333
$ws.Range("L131").Value = 5132.199900000001
$ws.Range("N131").Value = -15212.1999
$ws.Range("H137").Value = 1199.4
$ws.Range("I137").Value = 1199.4
$ws.Range("K137").Value = 3598.2
$ws.Range("M137").Value = 1501.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 26714
$ws.Range("J49").Value = 24499.5
$ws.Range("L49").Value = 24499.5
$ws.Range("N49").Value = -24867.5
$ws.Range("H126").Value = 102541.5
$ws.Range("J126").Value = 3249.5
$ws.Range("L126").Value = 9748.5
$ws.Range("N126").Value = -14688.5
$ws.Range("H132").Value = 1893.2
$ws.Range("I132").Value = 1991.5
$ws.Range("K132").Value = 5974.5
$ws.Range("M132").Value = -3444.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3367.158
$ws.Range("I16").Value = 3666.4666
$ws.Range("K16").Value = 3666.4666
$ws.Range("M16").Value = -3496.4666
$ws.Range("H40").Value = 4260.2354
$ws.Range("I40").Value = 3776.8125
$ws.Range("K40").Value = 3776.8125
$ws.Range("M40").Value = -3640.8125
$ws.Range("H42").Value = 19864
$ws.Range("I42").Value = 13682.667
$ws.Range("J42").Value = 24500
$ws.Range("K42").Value = 13682.667
$ws.Range("L42").Value = 24500
$ws.Range("M42").Value = -13119.667
$ws.Range("N42").Value = -25626
$ws.Range("H46").Value = 2081.5
$ws.Range("I46").Value = 815.8823
$ws.Range("J46").Value = 4472.1113
$ws.Range("K46").Value = 815.8823
$ws.Range("L46").Value = 4472.1113
$ws.Range("M46").Value = -627.8823
$ws.Range("N46").Value = -4848.1113
$ws.Range("H49").Value = 19864
$ws.Range("I49").Value = 13682.667
$ws.Range("J49").Value = 24500
$ws.Range("K49").Value = 13682.667
$ws.Range("L49").Value = 24500
$ws.Range("M49").Value = -13535.667
$ws.Range("N49").Value = -24794
$ws.Range("H93").Value = 2039.5
$ws.Range("I93").Value = 1639.4
$ws.Range("J93").Value = 2439.6
$ws.Range("K93").Value = 1639.4
$ws.Range("L93").Value = 2439.6
$ws.Range("M93").Value = -391.4000000000001
$ws.Range("N93").Value = -4935.6
$ws.Range("H94").Value = 60000
$ws.Range("J94").Value = 60000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352
$ws.Range("H100").Value = 4052
$ws.Range("I100").Value = 3248
$ws.Range("J100").Value = 5315.4287
$ws.Range("K100").Value = 3248
$ws.Range("L100").Value = 5315.4287
$ws.Range("M100").Value = -2707
$ws.Range("N100").Value = -6397.4287
$ws.Range("H132").Value = 7737.5
$ws.Range("I132").Value = 3874.1667
$ws.Range("J132").Value = 25122.5
$ws.Range("K132").Value = 11622.5001
$ws.Range("L132").Value = 75367.5
$ws.Range("M132").Value = -9092.500100000001
$ws.Range("N132").Value = -80427.5
$ws.Range("H136").Value = 2788
$ws.Range("I136").Value = 2744.5173
$ws.Range("J136").Value = 2998.1667
$ws.Range("K136").Value = 8233.5519
$ws.Range("L136").Value = 8994.500100000001
$ws.Range("M136").Value = -5683.5519
$ws.Range("N136").Value = -14094.5001
$ws.Range("H141").Value = 80715
$ws.Range("J141").Value = 80715
$ws.Range("L141").Value = 80715
$ws.Range("N141").Value = -91075

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 4
$ws.Range("K17").Value = 4
$ws.Range("M17").Value = 168
$ws.Range("H132").Value = 1725.1025
$ws.Range("I132").Value = 1510.0303
$ws.Range("K132").Value = 4530.090899999999
$ws.Range("M132").Value = -2000.090899999999
$ws.Range("H136").Value = 10992.405
$ws.Range("I136").Value = 9764.619000000001
$ws.Range("J136").Value = 12603.875
$ws.Range("K136").Value = 29293.857
$ws.Range("L136").Value = 37811.625
$ws.Range("M136").Value = -26743.857
$ws.Range("N136").Value = -42911.625
